$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Toggle the "Started" (Yes/No) flag in column C for the rows that changed.
$ws.Range("C3").Value  = "Yes"
$ws.Range("C9").Value  = "No"
$ws.Range("C31").Value = "Yes"
$ws.Range("C32").Value = "No"
$ws.Range("C42").Value = "No"
$ws.Range("C45").Value = "Yes"
$ws.Range("C51").Value = "No"
$ws.Range("C54").Value = "Yes"
$ws.Range("C55").Value = "Yes"
$ws.Range("C56").Value = "No"
$ws.Range("C76").Value = "Yes"
$ws.Range("C77").Value = "No"

# Re-create the view state: frozen header row stays frozen, scrolled so
# row 63 is the first visible row below the freeze, with C84 selected.
$ws.Activate() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C84").Select() | Out-Null
